$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Icam1", "Spn", "ECs", 3, 1, 144.843106, 434.529318, 0.6517202749316883, 0.6517202749316884, 2, 0.6666666666666666, 0.7483270000000001, 2.244981, 0.2123245966460564, 0.2123245966460564, 108.390006983662, 975.5100628529581, 0.1383762445009277, 0.1383762445009277),
    @("ECs", "Icam1", "Spn", "M2", 3, 1, 144.843106, 434.529318, 0.6517202749316883, 0.6517202749316884, 3, 1, 2.776121, 8.328363, 0.7876754033539436, 0.7876754033539437, 402.101988271826, 3618.917894446434, 0.5133440304307606, 0.5133440304307607),
    @("FAPs", "Icam1", "Spn", "ECs", 3, 1, 33.793597, 101.380791, 0.1520539909422453, 0.1520539909422453, 2, 0.6666666666666666, 0.7483270000000001, 2.244981, 0.2123245966460564, 0.2123245966460564, 25.288661062219, 227.597949559971, 0.03228480229523534, 0.03228480229523534),
    @("FAPs", "Icam1", "Spn", "M2", 3, 1, 33.793597, 101.380791, 0.1520539909422453, 0.1520539909422453, 3, 1, 2.776121, 8.328363, 0.7876754033539436, 0.7876754033539437, 93.81511429723699, 844.336028675133, 0.11976918864701, 0.11976918864701),
    @("M2", "Icam1", "Spn", "ECs", 3, 1, 41.33760066666667, 124.012802, 0.1859981687460937, 0.1859981687460937, 2, 0.6666666666666666, 0.7483270000000001, 2.244981, 0.2123245966460564, 0.2123245966460564, 30.93404269408467, 278.406384246762, 0.03949198615591949, 0.03949198615591949),
    @("M2", "Icam1", "Spn", "M2", 3, 1, 41.33760066666667, 124.012802, 0.1859981687460937, 0.1859981687460937, 3, 1, 2.776121, 8.328363, 0.7876754033539436, 0.7876754033539437, 114.7581813003473, 1032.823631703126, 0.1465061825901743, 0.1465061825901743),
    @("sCs", "Icam1", "Spn", "ECs", 3, 1, 2.273049333333333, 6.819148, 0.01022756537997252, 0.01022756537997253, 2, 0.6666666666666666, 0.7483270000000001, 2.244981, 0.2123245966460564, 0.2123245966460564, 1.700984188465333, 15.308857696188, 0.002171563693973837, 0.002171563693973837),
    @("sCs", "Icam1", "Spn", "M2", 3, 1, 2.273049333333333, 6.819148, 0.01022756537997252, 0.01022756537997253, 3, 1, 2.776121, 8.328363, 0.7876754033539436, 0.7876754033539437, 6.310259988302666, 56.792339894724, 0.008056001685998689, 0.00805600168599869)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
